# Remove the "SUM" helper column (J): it held a per-row =SUM(B:H) formula that is
# no longer needed now that the same condition is evaluated directly by the
# conditional-formatting rule below. Deleting the entire column also shifts the
# "Inclusion comments" column (L) one place to the left (becomes K).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").EntireColumn.Delete() | Out-Null

# Re-confirm the conditional formatting that strikes through any patient row whose
# scan count (B:H) is below 2 -- i.e. patients with only one scan.
$rng = $ws.Range("A2:I34")
$fc = $rng.FormatConditions.Item(1)
$fc.Font.Strikethrough = $true
$fc.Font.ThemeColor = 0
$fc.Font.TintAndShade = -0.34998626667073579

$ws.Range("K14").Select() | Out-Null
